$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new "Train r2" / "Test r2" values for the existing rows
$ws.Range("E4").Value = 0.536

$ws.Range("D5").Value = 0.5418
$ws.Range("E5").Value = 0.5375

$ws.Range("D6").Value = 0.6469
$ws.Range("E6").Value = 0.5502

$ws.Range("D7").Value = 0.7961
$ws.Range("E7").Value = 0.5969

$ws.Range("D8").Value = 0.756
$ws.Range("E8").Value = 0.5451

$ws.Range("D9").Value = 0.7376
$ws.Range("E9").Value = 0.5831

$ws.Range("D10").Value = 0.5914
$ws.Range("E10").Value = 0.5609

# New small italic "footer"-style cell far below the table (row 17)
$ws.Range("C17").Font.Italic = $true
$ws.Range("C17").Font.Size = 5
$ws.Range("C17").Font.Name = "Consolas"
$ws.Range("C17").Font.Family = 3
$ws.Range("C17").Font.Color = 9340031
$ws.Range("C17").VerticalAlignment = -4108

# Update the selection to match the diff
$ws.Range("E5").Select()
